$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before column D (old D:K -> new F:M), pushing the
# existing eight quarters two slots to the right and making room for two
# new, more-recent quarters (new D = most recent quarter, new E = next).
$ws.Range("D1:E1").EntireColumn.Insert()

# The insert copies the formatting of the column to its right onto the new
# blank columns is not reliable, so explicitly copy number formats/styles
# from the (now-shifted) old "D" column - which landed in F - onto the new
# D:E columns for the full data range.
$ws.Range("F5:F102").Copy()
$ws.Range("D5:E102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New quarter data for column D (most recent quarter) and column E (prior
# quarter), keyed by row number.
$newData = @(
    @{Row=7; D=43465; E=43373},
    @{Row=8; D=164300; E=169300},
    @{Row=9; D=44200; E=40000},
    @{Row=10; D=120100; E=129300},
    @{Row=12; D="NA"; E="NA"},
    @{Row=13; D=0; E=0},
    @{Row=14; D=0; E=0},
    @{Row=15; D=26700; E=26500},
    @{Row=17; D=140500; E=141000},
    @{Row=18; D=23800; E=28300},
    @{Row=20; D=700; E=100},
    @{Row=21; D=51200; E=54900},
    @{Row=22; D=7400; E=7000},
    @{Row=23; D=17100; E=21400},
    @{Row=24; D=7700; E=5600},
    @{Row=25; D=0; E=0},
    @{Row=26; D=9400; E=15800},
    @{Row=27; D=9400; E=15800},
    @{Row=28; D=0; E=0},
    @{Row=29; D=0; E="NA"},
    @{Row=30; D=0; E=0},
    @{Row=31; D=0; E=0},
    @{Row=32; D=-700; E=-100},
    @{Row=33; D=9400; E=15800},
    @{Row=34; D=0; E=0},
    @{Row=35; D=9400; E=15800},
    @{Row=38; D=43465; E=43373},
    @{Row=41; D=25500; E=17800},
    @{Row=42; D=0; E=0},
    @{Row=43; D=108900; E=107700},
    @{Row=44; D=0; E=0},
    @{Row=45; D=19600; E=21200},
    @{Row=46; D=153900; E=146700},
    @{Row=47; D=10300; E=10500},
    @{Row=48; D=41500; E=40900},
    @{Row=49; D=2394900; E=2417800},
    @{Row=50; D=0; E=0},
    @{Row=51; D=0; E=0},
    @{Row=52; D=0; E=0},
    @{Row=53; D=0; E=0},
    @{Row=54; D=2600500; E=2615800},
    @{Row=57; D=11600; E=8000},
    @{Row=58; D=26900; E=24000},
    @{Row=59; D=72200; E=81800},
    @{Row=60; D=110700; E=113900},
    @{Row=61; D=665300; E=678400},
    @{Row=62; D=197600; E=188700},
    @{Row=63; D=0; E=0},
    @{Row=64; D=0; E=0},
    @{Row=65; D=0; E=0},
    @{Row=66; D=973600; E=981000},
    @{Row=68; D=0; E=0},
    @{Row=69; D=0; E=0},
    @{Row=70; D=0; E=0},
    @{Row=71; D=0; E=0},
    @{Row=72; D=118200; E=128900},
    @{Row=73; D=0; E=0},
    @{Row=74; D=0; E=0},
    @{Row=75; D=0; E=0},
    @{Row=76; D=1626900; E=1634800},
    @{Row=77; D=0; E=0},
    @{Row=80; D=43465; E=43373},
    @{Row=81; D=9400; E=15800},
    @{Row=83; D=26700; E=26500},
    @{Row=84; D=0; E=0},
    @{Row=85; D=0; E=0},
    @{Row=86; D=0; E=0},
    @{Row=87; D=0; E=0},
    @{Row=88; D=0; E=0},
    @{Row=89; D=42500; E=50500},
    @{Row=91; D=-4300; E=-3500},
    @{Row=92; D=0; E=0},
    @{Row=93; D=0; E=0},
    @{Row=94; D=-4300; E=-3700},
    @{Row=96; D=0; E=0},
    @{Row=97; D=0; E=0},
    @{Row=98; D=0; E=0},
    @{Row=99; D=0; E=0},
    @{Row=100; D=-30600; E=-47300},
    @{Row=101; D=0; E=0},
    @{Row=102; D=7700; E=-600}

)

foreach ($item in $newData) {
    $ws.Cells.Item($item.Row, 4).Value2 = $item.D
    $ws.Cells.Item($item.Row, 5).Value2 = $item.E
}
